$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data (and row realignment for B15:E51)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.973.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.343.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.75%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.10"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.54"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0800"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.57"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.422.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.17%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.923.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.17%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.59"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +10.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.18"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.38"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.84"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.82%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.62"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0726"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.98"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.87"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.023.83"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.07"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.23%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.97"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.01"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.567.44"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.71"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.87%  "
